$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2:K11").Value = "h-"
